$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 79 (shifts old rows 79-119 down to 80-120)
$ws.Rows.Item(79).Insert()

# Fill the new row 79 with its data (same shape as the row below, new values)
$ws.Range("A79").Value = 11
$ws.Range("B79").Value = "Vega Monumental Concepción"
$ws.Range("C79").Value = "Bíobío"
$ws.Range("D79").Value = 44825
$ws.Range("E79").Value = 8
$ws.Range("F79").Value = 100112001
$ws.Range("G79").Value = "Berenjena"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 100
$ws.Range("K79").Value = 12000
$ws.Range("L79").Value = 13000
$ws.Range("M79").Value = 12500
$ws.Range("N79").Value = "$/caja 60 unidades"
$ws.Range("O79").Value = "Provincia de Limarí"
$ws.Range("P79").Value = 208
$ws.Range("Q79").Value = 60
$ws.Range("R79").Value = "Hortaliza"
